# Correct Ak from ELOW (values below 5 = 0)
#
# The upstream analysis recomputes biomass aggregates after flooring the
# "Ak" quantity derived from ELOW pH tiles to 0 whenever it was below 5.
# That correction removes a small (3.12) contribution from the Biomass_cal
# (F) and Biomass_fil (H) figures for ELOW/tile_29/T1,T2,T3 (rows 59-61),
# which in turn changes Biomass_overall (E = H + G, i.e.
# Biomass_fil + Biomass_npp) and the normalized Biomass_std_* columns
# (I, J, L), which are ratios against the ELOW/T0 baseline row (row 58).
# A few other rows (21, 37, 53) already had F = 0, and their
# Biomass_std_cal (J) values were pinned at a tiny non-zero floor (1e-4);
# those are now corrected to 0 as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple floor corrections: 1.0E-4 -> 0.0 -------------------------------
$ws.Range("J21").Value = 0.0
$ws.Range("J37").Value = 0.0
$ws.Range("J53").Value = 0.0

# --- Baseline row (ELOW, tile_29, T0) used to normalize std columns -------
$E58 = $ws.Range("E58").Value2
$F58 = $ws.Range("F58").Value2
$G58 = $ws.Range("G58").Value2
$H58 = $ws.Range("H58").Value2

# --- Amount to remove from Biomass_cal / Biomass_fil (and hence overall) --
$Ak = 3.12

# Rows 59, 60, 61 (ELOW, tile_29, T1/T2/T3): subtract the now-zeroed Ak
# contribution from F (Biomass_cal) and H (Biomass_fil), recompute
# E (Biomass_overall = H + G), and recompute the std columns as ratios
# to the baseline row 58.
foreach ($r in 59..61) {
    $G = $ws.Range("G$r").Value2

    $newF = $ws.Range("F$r").Value2 - $Ak
    $newH = $ws.Range("H$r").Value2 - $Ak
    $newE = $newH + $G

    $ws.Range("E$r").Value = $newE
    $ws.Range("F$r").Value = $newF
    $ws.Range("H$r").Value = $newH

    $ws.Range("I$r").Value = $newE / $E58
    $ws.Range("J$r").Value = $newF / $F58
    $ws.Range("L$r").Value = $newH / $H58
}
